$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 485.36
$ws.Range("I28").Value = 454.73914
$ws.Range("K28").Value = 454.73914
$ws.Range("M28").Value = 30.26085999999998
# Row 33
$ws.Range("H33").Value = 130.125
$ws.Range("I33").Value = 135.85715
$ws.Range("K33").Value = 135.85715
$ws.Range("M33").Value = 93.14285000000001
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
# Row 86
$ws.Range("H86").Value = 6889
$ws.Range("I86").Value = 2999
$ws.Range("J86").Value = 8185.6665
$ws.Range("K86").Value = 2999
$ws.Range("L86").Value = 8185.6665
$ws.Range("M86").Value = -1876
$ws.Range("N86").Value = -10431.6665
# Row 89
$ws.Range("H89").Value = 6889
$ws.Range("I89").Value = 2999
$ws.Range("J89").Value = 8185.6665
$ws.Range("K89").Value = 14995
$ws.Range("L89").Value = 40928.3325
$ws.Range("M89").Value = -9379
$ws.Range("N89").Value = -52160.3325
# Row 129
$ws.Range("H129").Value = 1839.4
$ws.Range("J129").Value = 2000
$ws.Range("L129").Value = 6000
$ws.Range("N129").Value = -16000
# Row 132
$ws.Range("H132").Value = 754.25
$ws.Range("I132").Value = 754.25
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2262.75
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 267.25
$ws.Range("N132").ClearContents()
# Row 137
$ws.Range("H137").Value = 1700
$ws.Range("I137").Value = 1700
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5100
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -2550
$ws.Range("N137").ClearContents()
# Row 138
$ws.Range("H138").Value = 9124.6875
$ws.Range("I138").Value = 7599.8
$ws.Range("J138").Value = 9817.817999999999
$ws.Range("K138").Value = 22799.4
$ws.Range("L138").Value = 29453.454
$ws.Range("M138").Value = -17659.4
$ws.Range("N138").Value = -39733.454

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 4949.8
$ws.Range("J45").Value = 4949.8
$ws.Range("L45").Value = 4949.8
$ws.Range("N45").Value = -5703.8
# Row 132
$ws.Range("H132").Value = 2830.6667
$ws.Range("I132").Value = 2830.6667
$ws.Range("K132").Value = 8492.000100000001
$ws.Range("M132").Value = -5962.000100000001
# Row 134
$ws.Range("H134").Value = 78500
$ws.Range("J134").Value = 78500
$ws.Range("L134").Value = 78500
$ws.Range("N134").Value = -88640

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
# Row 20
$ws.Range("H20").Value = 1798.8334
$ws.Range("I20").Value = 1798.8334
$ws.Range("K20").Value = 1798.8334
$ws.Range("M20").Value = -1551.8334

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 14
$ws.Range("H14").Value = 250
$ws.Range("J14").Value = 250
$ws.Range("L14").Value = 250
$ws.Range("N14").Value = -590
# Row 17
$ws.Range("H17").Value = 10333.333
# Row 31
$ws.Range("H31").Value = 2777.5454
$ws.Range("I31").Value = 2172.6667
$ws.Range("K31").Value = 2172.6667
$ws.Range("M31").Value = -1877.6667
# Row 34
$ws.Range("H34").Value = 2777.5454
$ws.Range("I34").Value = 2172.6667
$ws.Range("K34").Value = 2172.6667
$ws.Range("M34").Value = -1970.6667
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
# Row 86
$ws.Range("H86").Value = 10000
$ws.Range("I86").Value = 10000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 10000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -8877
$ws.Range("N86").ClearContents()
# Row 89
$ws.Range("H89").Value = 10000
$ws.Range("I89").Value = 10000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 50000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -44384
$ws.Range("N89").ClearContents()
# Row 107
$ws.Range("H107").Value = 637
$ws.Range("I107").Value = 664.4
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 664.4
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1255.6
$ws.Range("N107").Value = -4340
# Row 122
$ws.Range("H122").Value = 1610.4
$ws.Range("I122").Value = 1412
$ws.Range("J122").Value = 1742.6666
$ws.Range("K122").Value = 4236
$ws.Range("L122").Value = 5227.9998
$ws.Range("M122").Value = -1786
$ws.Range("N122").Value = -10127.9998

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 400000160
$ws.Range("I4").Value = 500000220
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 1500000660
$ws.Range("L4").Value = 12
$ws.Range("M4").Value = -1500000548
$ws.Range("N4").Value = -236
# Row 7
$ws.Range("H7").Value = 114.333336
$ws.Range("I7").Value = 165.5
$ws.Range("J7").Value = 12
$ws.Range("K7").Value = 496.5
$ws.Range("L7").Value = 36
$ws.Range("M7").Value = -384.5
$ws.Range("N7").Value = -260
# Row 9
$ws.Range("H9").Value = 425
$ws.Range("J9").Value = 400
$ws.Range("L9").Value = 1200
$ws.Range("N9").Value = -1648
# Row 15
$ws.Range("H15").Value = 1474.75
$ws.Range("J15").Value = 1474.75
$ws.Range("L15").Value = 4424.25
$ws.Range("N15").Value = -4704.25
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()
# Row 131
$ws.Range("H131").Value = 900
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
# Row 137
$ws.Range("H137").Value = 3298
$ws.Range("J137").Value = 3494.6667
$ws.Range("L137").Value = 10484.0001
$ws.Range("N137").Value = -20684.0001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 135
$ws.Range("H135").Value = 30000000
$ws.Range("J135").Value = 30000000
$ws.Range("L135").Value = 30000000
$ws.Range("N135").Value = -30010140

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 6299.3335
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
# Row 22
$ws.Range("H22").Value = 3778.4211
$ws.Range("I22").Value = 3069.5715
$ws.Range("K22").Value = 3069.5715
$ws.Range("M22").Value = -2774.5715
# Row 27
$ws.Range("H27").Value = 3778.4211
$ws.Range("I27").Value = 3069.5715
$ws.Range("K27").Value = 3069.5715
$ws.Range("M27").Value = -2962.5715
# Row 40
$ws.Range("H40").Value = 2374.25
$ws.Range("I40").Value = 2374.25
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2374.25
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2238.25
$ws.Range("N40").ClearContents()
# Row 43
$ws.Range("H43").Value = 13666.667
$ws.Range("J43").Value = 12000
$ws.Range("L43").Value = 12000
$ws.Range("N43").Value = -12386
# Row 56
$ws.Range("H56").Value = 30000
$ws.Range("I56").Value = 30000
$ws.Range("K56").Value = 30000
$ws.Range("M56").Value = -29309
# Row 105
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
# Row 122
$ws.Range("H122").Value = 1954
$ws.Range("I122").Value = 1953.5
$ws.Range("J122").Value = 1956
$ws.Range("K122").Value = 5860.5
$ws.Range("L122").Value = 5868
$ws.Range("M122").Value = -3410.5
$ws.Range("N122").Value = -10768

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Range("H74").Value = 37500
$ws.Range("I74").Value = 35000
$ws.Range("K74").Value = 35000
$ws.Range("M74").Value = -34064
# Row 77
$ws.Range("H77").Value = 37500
$ws.Range("I77").Value = 35000
$ws.Range("K77").Value = 105000
$ws.Range("M77").Value = -100320
# Row 95
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492
# Row 132
$ws.Range("H132").Value = 2899.75
$ws.Range("J132").Value = 2866.3333
$ws.Range("L132").Value = 8598.999899999999
$ws.Range("N132").Value = -13658.9999
# Row 136
$ws.Range("H136").Value = 1135
$ws.Range("I136").Value = 1135.6
$ws.Range("K136").Value = 3406.8
$ws.Range("M136").Value = -856.7999999999997

